# Applies the "Updated cryptos list" price/volume refresh to Sheet1 (rows 2-51).
# D column = Price, E column = Volume(1h) change percentage.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Scratch cell used to force plain-numeric-looking price strings (e.g. "6.20")
# to be written as literal text, matching the source data's inlineStr cells,
# instead of being auto-coerced into numeric cells by COM type inference.
$scratch = $ws.Range("Z1")

function Set-TextValue($range, $text) {
    $scratch.NumberFormat = "@"
    $scratch.Value = $text
    $scratch.Copy()
    $range.PasteSpecial(-4163)
}

$ws.Range("D2").Value = "66.877.80"
$ws.Range("E2").Value = "  -1.00%  "

$ws.Range("D3").Value = "3.459.85"
$ws.Range("E3").Value = "  -2.15%  "

$ws.Range("E4").Value = "  +0.04%  "

Set-TextValue $ws.Range("D5") "591.36"
$ws.Range("E5").Value = "  -1.08%  "

Set-TextValue $ws.Range("D6") "174.89"
$ws.Range("E6").Value = "  +0.54%  "

$ws.Range("E7").Value = "  +0.07%  "

Set-TextValue $ws.Range("D8") "0.584"
$ws.Range("E8").Value = "  -1.59%  "

Set-TextValue $ws.Range("D9") "0.128"
$ws.Range("E9").Value = "  -4.76%  "

Set-TextValue $ws.Range("D10") "7.06"
$ws.Range("E10").Value = "  -3.57%  "

Set-TextValue $ws.Range("D11") "0.423"
$ws.Range("E11").Value = "  -3.07%  "

$ws.Range("D12").Value = "4.066.36"
$ws.Range("E12").Value = "  -1.81%  "

Set-TextValue $ws.Range("D13") "30.57"
$ws.Range("E13").Value = "  +5.95%  "

$ws.Range("E14").Value = "  -0.23%  "

$ws.Range("D15").Value = "66.910.03"
$ws.Range("E15").Value = "  -0.77%  "

Set-TextValue $ws.Range("D16") "0.0000174"
$ws.Range("E16").Value = "  -4.59%  "

$ws.Range("D17").Value = "3.468.42"
$ws.Range("E17").Value = "  -1.88%  "

Set-TextValue $ws.Range("D18") "6.20"
$ws.Range("E18").Value = "  -2.64%  "

Set-TextValue $ws.Range("D19") "14.25"
$ws.Range("E19").Value = "  +0.22%  "

Set-TextValue $ws.Range("D20") "385.06"
$ws.Range("E20").Value = "  -3.22%  "

Set-TextValue $ws.Range("D21") "7.81"
$ws.Range("E21").Value = "  -2.45%  "

Set-TextValue $ws.Range("D22") "72.38"
$ws.Range("E22").Value = "  -1.68%  "

Set-TextValue $ws.Range("D23") "0.996"
$ws.Range("E23").Value = "  -0.37%  "

$ws.Range("E24").Value = "  -0.26%  "

Set-TextValue $ws.Range("D25") "0.530"
$ws.Range("E25").Value = "  -1.98%  "

Set-TextValue $ws.Range("D26") "0.0000120"
$ws.Range("E26").Value = "  -3.16%  "

Set-TextValue $ws.Range("D27") "10.23"
$ws.Range("E27").Value = "  -0.62%  "

$ws.Range("E28").Value = "  -2.57%  "

Set-TextValue $ws.Range("D29") "0.995"
$ws.Range("E29").Value = "  -0.31%  "

Set-TextValue $ws.Range("D30") "6.06"
$ws.Range("E30").Value = "  -4.05%  "

Set-TextValue $ws.Range("D31") "1.42"
$ws.Range("E31").Value = "  -4.07%  "

$ws.Range("E32").Value = "  -2.86%  "

Set-TextValue $ws.Range("D33") "23.32"
$ws.Range("E33").Value = "  -3.35%  "

Set-TextValue $ws.Range("D34") "7.23"
$ws.Range("E34").Value = "  -2.36%  "

Set-TextValue $ws.Range("D35") "1.60"
$ws.Range("E35").Value = "  -2.22%  "

Set-TextValue $ws.Range("D36") "163.05"
$ws.Range("E36").Value = "  -0.54%  "

Set-TextValue $ws.Range("D37") "0.868"
$ws.Range("E37").Value = "  -3.35%  "

Set-TextValue $ws.Range("D38") "1.90"
$ws.Range("E38").Value = "  -1.16%  "

Set-TextValue $ws.Range("D39") "6.99"
$ws.Range("E39").Value = "  -0.09%  "

Set-TextValue $ws.Range("D40") "27.15"
$ws.Range("E40").Value = "  -0.74%  "

Set-TextValue $ws.Range("D41") "4.58"
$ws.Range("E41").Value = "  -3.43%  "

Set-TextValue $ws.Range("D42") "26.13"
$ws.Range("E42").Value = "  -1.83%  "

$ws.Range("D43").Value = "2.774.46"
$ws.Range("E43").Value = "  -1.49%  "

Set-TextValue $ws.Range("D44") "0.0717"
$ws.Range("E44").Value = "  -4.47%  "

Set-TextValue $ws.Range("D45") "2.56"
$ws.Range("E45").Value = "  -3.62%  "

Set-TextValue $ws.Range("D46") "42.07"
$ws.Range("E46").Value = "  -2.22%  "

Set-TextValue $ws.Range("D47") "0.0296"
$ws.Range("E47").Value = "  -4.85%  "

Set-TextValue $ws.Range("D48") "336.65"
$ws.Range("E48").Value = "  -1.77%  "

$ws.Range("E49").Value = "  -3.67%  "

Set-TextValue $ws.Range("D50") "33.02"
$ws.Range("E50").Value = "  -2.84%  "

Set-TextValue $ws.Range("D51") "6.33"
$ws.Range("E51").Value = "  -3.26%  "

# Clean up the scratch cell so it does not leave stray content/formatting behind.
$scratch.Clear()